# Zeiterfassung.xlsx update
# - Add a new time-tracking entry row to the "Julian" sheet:
#     Datum:      03.07.2022
#     Taetigkeit: Frontend Design Login- und Landing-Page
#     Aufwand:    1.4
#   The row is inserted right above the existing "Gesamt" (total) row,
#   which shifts down by one and whose SUM formula grows to include the
#   newly inserted row.
# - Make "Julian" the active sheet/tab, with cell J5 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Julian")

# Insert a new row above the current "Gesamt" row (row 5), pushing it
# down to row 6. This also shifts the merged A5:B5 range and the SUM
# formula's relative row references along with it.
$ws.Range("A5:C5").EntireRow.Insert()

# --- Datum (A5) ---------------------------------------------------------
# The date column stores plain text labels (e.g. "26.06.2022"), not real
# date values. Assigning the literal string directly would make Excel
# auto-parse it into a date serial number, so build it as a text formula
# first (never auto-parsed) and then collapse it down to a plain value.
$ws.Range("A5").Formula = '="03.07.2022"'
$ws.Range("A5").Copy()
$ws.Range("A5").PasteSpecial(-4163) | Out-Null  # xlPasteValues
# Re-apply the same formatting the other data rows use (vertical-center,
# no explicit number format) by copying it over from row 3.
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Taetigkeit (B5) -----------------------------------------------------
$ws.Range("B5").Value = "Frontend Design Login- und Landing-Page"
$ws.Range("B3").Copy()
$ws.Range("B5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Aufwand (C5) --------------------------------------------------------
$ws.Range("C5").Value = 1.4
$ws.Range("C3").Copy()
$ws.Range("C5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# The new row's wrapped Taetigkeit text needs the same extra row height
# as the other multi-line entries.
$ws.Rows.Item(5).RowHeight = 28.8

# Extend the Gesamt row's SUM formula to cover the newly inserted row.
$ws.Range("C6").Formula = "=SUM(C3:C5)"

# Make "Julian" the active sheet and set its new selection.
$ws.Activate()
$ws.Range("J5").Select()
